$d = $word.ActiveDocument

# Remove the _GoBack bookmark that splits the bold run "R" | "TF and TXT"
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# Merge the split text "DOCX, DOC, PDF, HTML, XPS, R" + "TF and TXT"
# into a single continuous run "DOCX, DOC, PDF, HTML, XPS, RTF and TXT"
$d.Content.Find.Execute("DOCX, DOC, PDF, HTML, XPS, R" + "TF and TXT", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DOCX, DOC, PDF, HTML, XPS, RTF and TXT", 2)
